$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 2 (shifts existing rows 2..49 down to 3..50)
$ws.Rows.Item(2).Insert()

# The insert operation copies formatting from the row above (the bold
# header). Clear that so the new row matches the plain data-row look,
# then re-apply just the date number format (style used by column D)
# by copying the format from the date cell directly below.
$ws.Range("A2:R2").ClearFormats()
$ws.Range("D3").Copy()
$ws.Range("D2").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Fill in the new row's values (week of 2021-11-15, Macroferia Regional
# de Talca, Maule / Provincia de Linares, Esparragos Verde Primera)
$ws.Range("A2").Value = 5
$ws.Range("B2").Value = "Macroferia Regional de Talca"
$ws.Range("C2").Value = "Maule"
$ws.Range("D2").Value = "2021-11-15"
$ws.Range("E2").Value = 7
$ws.Range("F2").Value = 300000000
$ws.Range("G2").Value = "Espárragos"
$ws.Range("H2").Value = "Verde"
$ws.Range("I2").Value = "Primera"
$ws.Range("J2").Value = 3000
$ws.Range("K2").Value = 1200
$ws.Range("L2").Value = 1200
$ws.Range("M2").Value = 1200
$ws.Range("N2").Value = "$/kilo"
$ws.Range("O2").Value = "Provincia de Linares"
$ws.Range("P2").Value = 1200
$ws.Range("Q2").Value = 1
$ws.Range("R2").Value = "Hortaliza"
